# Daily attendance processing - 2025-10-31 02:59:21
# Normalize the "Recorded By" (column G) lists so the names/emails are
# sorted into ascending (ordinal / case-sensitive, ASCII) order within
# each cell. Entries are separated by ", ".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Value2

    if ($text -eq $null) { continue }
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Count -lt 2) { continue }

    # Build an ordinal (ASCII, case-sensitive) sort key for each entry so
    # that "System" (capital S = 83) sorts before lowercase entries such
    # as "system"/"backup@..." (>= 97), matching .NET's CompareOrdinal
    # rather than the culture-aware default comparer. A leading non-digit
    # prefix keeps Sort-Object from treating the zero-padded numeric key
    # as an actual number.
    $keyed = @()
    foreach ($part in $parts) {
        $key = "k"
        foreach ($ch in $part.ToCharArray()) {
            $code = [int][char]$ch
            $padded = ([string]$code).PadLeft(6, '0')
            $key = "$key$padded"
        }
        $obj = [PSCustomObject]@{ Orig = $part; Key = $key }
        $keyed = $keyed + @($obj)
    }

    $sortedObjs = $keyed | Sort-Object -Property Key
    $sorted = @()
    foreach ($o in $sortedObjs) {
        $sorted = $sorted + @($o.Orig)
    }

    $newText = $sorted -join ", "

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
